# Saldo.xlsx update
# 1) Add a new account row (Elizabeth) right above the existing first data row.
# 2) Remove the Ana and Elaine account rows.
# 3) Add a new account row (Andrea) right after the Camila account row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$text) {
    # Force the cell to keep a leading-zero account number as literal text
    # (rather than being auto-coerced to a number), then drop the
    # temporary "@" text format so the cell keeps the sheet's default style.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- 1) Insert the Elizabeth row before row 2 (004359408 / Hepta) ---
$ws.Rows.Item(2).Insert()
Set-TextCell $ws.Range("A2") "008297098"
Set-TextCell $ws.Range("B2") "Elizabeth"
$ws.Range("C2").Value = 447196.24

# --- 2) Delete the Ana (004267119) and Elaine (008115265) rows ---
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# --- 3) Insert the Andrea row right after the Camila row (366.69) ---
$ws.Rows.Item(17).Insert()
Set-TextCell $ws.Range("A17") "005186167"
Set-TextCell $ws.Range("B17") "Andrea"
$ws.Range("C17").Value = 349.98
